# Auto update Excel log
# Appends new sensor/event log rows to several worksheets, as produced by
# the SeniorConnect monitoring system on 2026-02-01 around 14:33-14:36.

$wb = $excel.ActiveWorkbook

function Append-Rows {
    param($SheetName, $StartRow, $Rows)

    $ws = $wb.Worksheets.Item($SheetName)
    $rowCount = $Rows.Count
    $endRow = $StartRow + $rowCount - 1

    # Force the new cells to be plain text so date/time-looking values
    # (e.g. "2026-02-01", "14:35:39") are not auto-converted into Excel
    # date/time serial numbers. Using the "A1:F2" colon-range syntax here
    # (rather than the two-argument Range(start, end) form) is required
    # for the format to actually apply across the whole block.
    $targetRange = $ws.Range("A" + $StartRow + ":F" + $endRow)
    $targetRange.NumberFormat = "@"

    for ($i = 0; $i -lt $rowCount; $i++) {
        $r = $StartRow + $i
        $rowData = $Rows[$i]
        $ws.Range("A$r").Value = $rowData[0]
        $ws.Range("B$r").Value = $rowData[1]
        $ws.Range("C$r").Value = $rowData[2]
        $ws.Range("D$r").Value = $rowData[3]
        $ws.Range("E$r").Value = $rowData[4]
        $ws.Range("F$r").Value = $rowData[5]
    }
}

# ALERTS sheet: three new CRITICAL fall-detection alerts (rows 5-7)
Append-Rows "ALERTS" 5 @(
    ,@("2026-02-01", "14:35:39", "14:00", "Living Room", "CRITICAL", "FALL_DETECTED")
    ,@("2026-02-01", "14:35:42", "14:00", "Living Room", "CRITICAL", "FALL_DETECTED")
    ,@("2026-02-01", "14:35:45", "14:00", "Living Room", "CRITICAL", "FALL_DETECTED")
)

# Proximity sheet: Living Room Main Door enter/exit events (rows 24-25)
Append-Rows "Proximity" 24 @(
    ,@("2026-02-01", "14:35:45", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
    ,@("2026-02-01", "14:35:49", "14:00", "Living Room Main Door", "EXIT", "User EXITED Living Room Main Door")
)

# mmWave sheet: additional presence-detected readings (rows 4-8)
Append-Rows "mmWave" 4 @(
    ,@("2026-02-01", "14:33:55", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")
    ,@("2026-02-01", "14:35:45", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")
    ,@("2026-02-01", "14:35:49", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")
    ,@("2026-02-01", "14:35:58", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")
    ,@("2026-02-01", "14:36:09", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

# Camera sheet: one additional image-captured event (row 14)
Append-Rows "Camera" 14 @(
    ,@("2026-02-01", "14:35:49", "14:00", "Living Room Main Door", "Image Captured", "Active")
)
